$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write an account number so it is stored as text with leading
# zeros preserved (like the existing "Conta" column), without leaving a
# stray quote-prefix style behind.
function Set-Conta($rowNum, $conta) {
    $cell = $ws.Cells.Item($rowNum, 1)
    $cell.Value = "'" + $conta
    $cell.ClearFormats()
}

function Insert-DataRow($rowNum, $conta, $nome, $saldo) {
    $ws.Rows.Item($rowNum).Insert()
    Set-Conta $rowNum $conta
    $ws.Cells.Item($rowNum, 2).Value = $nome
    $ws.Cells.Item($rowNum, 3).Value = $saldo
}

function Delete-DataRow($rowNum) {
    $ws.Rows.Item($rowNum).Delete()
}

# Apply changes from the bottom of the sheet upward so that earlier
# (lower) row numbers referenced below remain valid while we work.

# Row 155: remove 001090818 MARIA 0.09
Delete-DataRow 155

# Row 73: remove 004242237 MARIAH 60.32
Delete-DataRow 73

# Rows 24-25: 004646727 RENATA 244.18 / 004498637 TIAGO 243.4
# -> becomes 004498637 TIAGO 243.4 / 004646727 RENATA 230.58
Delete-DataRow 25
$ws.Cells.Item(24, 3).Value = 230.58
Insert-DataRow 24 "004498637" "TIAGO" 243.4

# Row 14: remove 004474776 GILSON 973.01
Delete-DataRow 14

# Before row 10 (004467884 ANA 2612.52): insert 4 new rows
Insert-DataRow 10 "004452597" "LARA" 14481.31
Insert-DataRow 11 "004242237" "MARIAH" 8863.86
Insert-DataRow 12 "004575632" "ADELE" 5246.6
Insert-DataRow 13 "001090818" "MARIA" 4703.25

# Row 6: 004224405 MILA 17710.87 -> 002786022 PAULO 17904.43, plus a new
# 004381095 THIAGO 16051.42 row right after it.
Set-Conta 6 "002786022"
$ws.Cells.Item(6, 2).Value = "PAULO"
$ws.Cells.Item(6, 3).Value = 17904.43
Insert-DataRow 7 "004381095" "THIAGO" 16051.42

# Before row 5 (004364200 BLOCO 19918.99): insert new 004313254 GUSTAVO row
Insert-DataRow 5 "004313254" "GUSTAVO" 22633.57
